$d = $word.ActiveDocument

# --- 1. Update the date "25.02.2015" -> "10.04.2015" in the "Lodz, dnia ..." line ---
$r = $d.Content
$found = $r.Find.Execute("25.02.2015", $true, $false, $false, $false, $false, $true, 1, $false, "10.04.2015", 2)

# $r now covers the freshly written "10.04.2015" text.
$base = $r.Start

# The replace above tends to also fold the preceding "Lodz, dnia " run into the
# edited run. Re-establish that original run boundary (position 0, i.e. right
# before the "1" of "10") by touching it with a transient bookmark: adding and
# immediately deleting a bookmark at a text position forces the engine to keep
# the runs split there, without leaving any bookmark behind.
$pHeal = $base
$bmHeal = $d.Range($pHeal, $pHeal)
$d.Bookmarks.Add("TmpSplit0", $bmHeal)
$d.Bookmarks("TmpSplit0").Delete()

# Split "10" | ".04.2015" and drop the real "_GoBack" bookmark right there
# (between "10" and ".04"), exactly like Word leaves behind after the last
# edit position. Re-adding a bookmark named "_GoBack" moves the existing one
# from wherever it used to be (end of doc) to this new spot.
$p1 = $base + 2
$bm1 = $d.Range($p1, $p1)
$d.Bookmarks.Add("_GoBack", $bm1)

# Split ".04" | ".2015" using the same transient-bookmark trick.
$p2 = $base + 5
$bm2 = $d.Range($p2, $p2)
$d.Bookmarks.Add("TmpSplit2", $bm2)
$d.Bookmarks("TmpSplit2").Delete()

# Split ".2015" | " roku" using the same transient-bookmark trick.
$p3 = $base + 10
$bm3 = $d.Range($p3, $p3)
$d.Bookmarks.Add("TmpSplit3", $bm3)
$d.Bookmarks("TmpSplit3").Delete()
